# Insert a new "Match ID" column at the front of the sheet (shifts B:W -> C:X)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("A:A").Insert()

# Header label for the new column
$ws.Range("A1").Value = "Match ID"

# New column's cells share the bold "header/id" font used elsewhere (no border)
$ws.Range("A1:A19").Font.Bold = $true

# Match ID value (17) for every real data row, including the hidden
# rows 4-19 is the visible body; row 20 is the hidden totals row.
$ws.Range("A4:A19").Value = 17
$ws.Cells.Item(20, 1).Value = 17

# Undo the implicit row-height bump that writing into the hidden row causes
$ws.Rows.Item(20).AutoFit()

# Reflect the edit in the sheet's selection (A1:A19 = the column just filled)
$ws.Range("A1:A19").Select() | Out-Null
